# Generate Report for Handoff
#
# Files 4-7 ("0949e904...", "8bf1144a...", "b00a0896...", "fa4f4dc8...") in both
# the zh-cn and de-de localization-status sheets were "low" priority and are now
# being handed off, so their Priority flips to "ht" and their
# "Latest Handoff Datetime" is refreshed to the moment the handoff report was
# generated. The Overview sheet's "Latest HO Xliff Generate Date" column mirrors
# that same timestamp for those rows, so it is updated to match as well.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn: rows 4-7 -> Priority (col E) "low" -> "ht"; Latest Handoff Datetime (col H) refreshed
for ($r = 4; $r -le 7; $r++) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-14 00:41:40"
}

# de-de: rows 4-7 -> Priority (col E) "low" -> "ht"; Latest Handoff Datetime (col H) refreshed
for ($r = 4; $r -le 7; $r++) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-14 00:41:47"
}

# Overview: rows 4-7, col G (Latest HO Xliff Generate Date) mirrors the de-de
# handoff timestamp that was just refreshed above.
for ($r = 4; $r -le 7; $r++) {
    $overview.Cells.Item($r, 7).Value = "2016-08-14 00:41:47"
}
